# "fixed export and fixing maps"
#
# The sheet originally showed the municipality area for three census
# years (1989 / 2002 / 2014) plus a "(according to the population
# census data)" subtitle under the title. The edit simplifies the
# table down to just the current (2014) figure:
#   - remove the subtitle row under the title
#   - remove the 1989 and 2002 columns, keeping only the 2014 column
#   - give the remaining rows a bit more breathing room (taller rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 held "(according to the population census data)" - drop it,
# shifting everything below up by one row.
$ws.Rows.Item(2).Delete()

# Columns B and C held the 1989 and 2002 figures - drop both, shifting
# the 2014 column (originally D) left into column B.
$ws.Columns.Item(2).Delete()
$ws.Columns.Item(2).Delete()

# Give the (now 5) used rows a taller, more readable height.
$ws.Range("A1:B5").RowHeight = 20.1
